# multiple_inheritance.pptx — "slides: minor updates to sessions 5 and 6"
#
# 1) Refresh the cached "datetimeFigureOut" footer text (5/8/22 -> 5/18/22)
#    on the slide master, every slide layout, and the notes master.
# 2) Session 5 (slide 6) - reword the last bullet about `super`.
# 3) Session 6 (slide 13) - reword the last bullet about changing functionality.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            if ($shp.TextFrame.TextRange.Text -eq "5/8/22") {
                $shp.TextFrame.TextRange.Text = "5/18/22"
            }
        }
    }
}

# --- Slide master ---
Update-DatePlaceholder($p.SlideMaster)

# --- Every slide layout ---
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder($layouts.Item($L))
}

# --- Notes master ---
Update-DatePlaceholder($p.NotesMaster)

# --- Session 5 / slide 6: "... super <text>" bullet ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$run6 = $tr6.Characters(297, 45)
$run6.Text = " may not always do what you want"

# --- Session 6 / slide 13: "Changing functionality ..." bullet ---
$s13 = $p.Slides.Item(13)
$sh13 = $s13.Shapes.Item(2)
$tr13 = $sh13.TextFrame.TextRange
$run13 = $tr13.Characters(381, 85)
$run13.Text = "With this approach, changing functionality requires changing the lines that implement that functionality, which means it’s good code."
